$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells to text so values like "1.000" or "0.01512" are not
# auto-converted to numbers by Excel, matching the original inline-string formatting.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D27","D29","D30","D31","D32","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.592.82"
$ws.Range("E2").Value = "  +4.07%  "
$ws.Range("D3").Value = "1.743.17"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "246.69"
$ws.Range("E5").Value = "  +4.44%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.4824"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("D8").Value = "0.2697"
$ws.Range("E8").Value = "  +3.99%  "
$ws.Range("D9").Value = "0.06265"
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("D10").Value = "1.744.12"
$ws.Range("E10").Value = "  +4.49%  "
$ws.Range("D11").Value = "0.07127"
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("D12").Value = "15.88"
$ws.Range("E12").Value = "  +7.48%  "
$ws.Range("D13").Value = "0.6246"
$ws.Range("E13").Value = "  +6.83%  "
$ws.Range("D14").Value = "4.519"
$ws.Range("E14").Value = "  +3.64%  "
$ws.Range("D15").Value = "77.43"
$ws.Range("D16").Value = "0.9999"
$ws.Range("D17").Value = "26.598.01"
$ws.Range("E17").Value = "  +4.11%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000006903"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "11.80"
$ws.Range("E20").Value = "  +3.45%  "
$ws.Range("D21").Value = "1.968.97"
$ws.Range("E21").Value = "  +4.60%  "
$ws.Range("D22").Value = "4.616"
$ws.Range("E22").Value = "  +4.02%  "
$ws.Range("D23").Value = "8.883"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").Value = "5.369"
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").Value = "136.45"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("D27").Value = "1.817"
$ws.Range("E27").Value = "  +6.15%  "
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("D29").Value = "106.88"
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("D30").Value = "4.016"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "3.743"
$ws.Range("E31").Value = "  +3.19%  "
$ws.Range("D32").Value = "0.07890"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("E33").Value = "  +6.69%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.616"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.6407"
$ws.Range("E35").Value = "  +5.87%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  +5.01%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "0.9329"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("B38").Value = "Quant"
$ws.Range("C38").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D38").Value = "113.70"
$ws.Range("E38").Value = "  +13.33%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.438"
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("D40").Value = "1.991"
$ws.Range("E40").Value = "  +8.28%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "5.810"
$ws.Range("E41").Value = "  +18.37%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "0.9998"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.01512"
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.3918"
$ws.Range("E44").Value = "  +4.71%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "0.1216"
$ws.Range("E45").Value = "  +9.23%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "6.763"
$ws.Range("E46").Value = "  +9.23%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.05336"
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "7.954"
$ws.Range("E48").Value = "  +6.43%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "30.79"
$ws.Range("E49").Value = "  +3.43%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.264"
$ws.Range("E50").Value = "  +5.15%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.3456"
$ws.Range("E51").Value = "  +4.16%  "

# Remove the temporary text-format styling so the cells keep the default (unstyled) look
foreach ($ref in $priceCells) {
    $ws.Range($ref).ClearFormats()
}

Write-Host "Applied cryptos list update"